# Weekly update: insert a new row of data at row 3 (most recent week's
# "Especial" quality entry), pushing all prior rows (old 3..19) down to
# (4..20). The new row reuses the constant/descriptive columns from the
# old row 3 (same Mercado/Region/Producto/Categoria/Variedad/Calidad/
# Unidad/Origen/Kg-unidad) and carries fresh Fecha/Volumen/Precio/
# Precio-$/Kg figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the row-3 values we need to preserve before the insert shifts
# everything down.
$colA = $ws.Range("A3").Value2
$colB = $ws.Range("B3").Value2
$colC = $ws.Range("C3").Value2
$colE = $ws.Range("E3").Value2
$colF = $ws.Range("F3").Value2
$colG = $ws.Range("G3").Value2
$colH = $ws.Range("H3").Value2
$colI = $ws.Range("I3").Value2
$colJ = $ws.Range("J3").Value2
$colK = $ws.Range("K3").Value2
$colL = $ws.Range("L3").Value2
$colQ = $ws.Range("Q3").Value2
$colR = $ws.Range("R3").Value2
$colT = $ws.Range("T3").Value2

# Push old rows 3..19 down to 4..20 (new blank row 3 created).
$ws.Rows.Item(3).Insert()

# Fill the new row 3 with the latest week's entry.
$ws.Range("A3").Value = $colA
$ws.Range("B3").Value = $colB
$ws.Range("C3").Value = $colC
$ws.Range("D3").Value = 44649
$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat
$ws.Range("E3").Value = $colE
$ws.Range("F3").Value = $colF
$ws.Range("G3").Value = $colG
$ws.Range("H3").Value = $colH
$ws.Range("I3").Value = $colI
$ws.Range("J3").Value = $colJ
$ws.Range("K3").Value = $colK
$ws.Range("L3").Value = $colL
$ws.Range("M3").Value = 18
$ws.Range("N3").Value = 330000
$ws.Range("O3").Value = 340000
$ws.Range("P3").Value = 335000
$ws.Range("Q3").Value = $colQ
$ws.Range("R3").Value = $colR
$ws.Range("S3").Value = 744
$ws.Range("T3").Value = $colT
